$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "Переводимый язык (основной диалект) [система письма] @ источник данных | обработка"
$ws.Range("C1").Value = "Другой язык (основной диалект) [система письма] @  источник данных | обработка"

$ws.Range("C2").Select()
